# Apply updated crypto price/volume data to the worksheet.
# Values in column D are prefixed with a leading apostrophe so Excel stores
# them as text (matching the original inlineStr text cells) instead of
# auto-converting number-looking strings (e.g. "3.48") into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'65.946.99"
$ws.Range("E2").Value = "  +0.61%  "

# Row 3
$ws.Range("D3").Value = "'3.323.91"
$ws.Range("E3").Value = "  +1.53%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.33%  "

# Row 5
$ws.Range("D5").Value = "'188.73"
$ws.Range("E5").Value = "  +4.69%  "

# Row 6
$ws.Range("D6").Value = "'555.34"
$ws.Range("E6").Value = "  -0.23%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.581"
$ws.Range("E8").Value = "  -1.33%  "

# Row 9
$ws.Range("D9").Value = "'3.314.84"
$ws.Range("E9").Value = "  +1.73%  "

# Row 10
$ws.Range("E10").Value = "  -3.42%  "

# Row 11
$ws.Range("D11").Value = "'0.581"
$ws.Range("E11").Value = "  -1.07%  "

# Row 12
$ws.Range("D12").Value = "'46.33"
$ws.Range("E12").Value = "  -2.46%  "

# Row 13
$ws.Range("D13").Value = "'0.0000269"
$ws.Range("E13").Value = "  +1.68%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'8.57"
$ws.Range("E14").Value = "  +0.31%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "'3.857.36"
$ws.Range("E15").Value = "  +1.03%  "

# Row 16
$ws.Range("D16").Value = "'593.97"
$ws.Range("E16").Value = "  -6.28%  "

# Row 17
$ws.Range("D17").Value = "'65.940.69"
$ws.Range("E17").Value = "  +0.41%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.117"
$ws.Range("E18").Value = "  +0.86%  "

# Row 19
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'17.88"
$ws.Range("E19").Value = "  +0.21%  "

# Row 20
$ws.Range("D20").Value = "'3.329.32"
$ws.Range("E20").Value = "  +1.27%  "

# Row 21
$ws.Range("D21").Value = "'10.98"
$ws.Range("E21").Value = "  -3.37%  "

# Row 22
$ws.Range("D22").Value = "'0.897"
$ws.Range("E22").Value = "  -0.83%  "

# Row 23
$ws.Range("D23").Value = "'18.48"
$ws.Range("E23").Value = "  +3.58%  "

# Row 24
$ws.Range("D24").Value = "'5.01"
$ws.Range("E24").Value = "  +0.75%  "

# Row 25
$ws.Range("D25").Value = "'99.18"
$ws.Range("E25").Value = "  -7.25%  "

# Row 26
$ws.Range("D26").Value = "'3.95"
$ws.Range("E26").Value = "  -1.06%  "

# Row 27
$ws.Range("D27").Value = "'5.96"
$ws.Range("E27").Value = "  -0.77%  "

# Row 28
$ws.Range("E28").Value = "  +1.92%  "

# Row 29
$ws.Range("D29").Value = "'9.49"
$ws.Range("E29").Value = "  -0.14%  "

# Row 30
$ws.Range("E30").Value = "  -1.46%  "

# Row 31
$ws.Range("D31").Value = "'30.40"
$ws.Range("E31").Value = "  +0.30%  "

# Row 32
$ws.Range("D32").Value = "'6.69"
$ws.Range("E32").Value = "  +5.58%  "

# Row 33
$ws.Range("D33").Value = "'3.93"
$ws.Range("E33").Value = "  -0.91%  "

# Row 34
$ws.Range("D34").Value = "'583.24"
$ws.Range("E34").Value = "  +5.01%  "

# Row 35
$ws.Range("D35").Value = "'10.94"
$ws.Range("E35").Value = "  -1.07%  "

# Row 36
$ws.Range("E36").Value = "  -1.06%  "

# Row 37
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "'3.701.31"
$ws.Range("E38").Value = "  +1.18%  "

# Row 39
$ws.Range("D39").Value = "'56.14"
$ws.Range("E39").Value = "  -1.75%  "

# Row 40
$ws.Range("D40").Value = "'3.48"
$ws.Range("E40").Value = "  -8.77%  "

# Row 41
$ws.Range("D41").Value = "'33.60"
$ws.Range("E41").Value = "  +4.89%  "

# Row 42
$ws.Range("D42").Value = "'0.0₃0701"
$ws.Range("E42").Value = "  -1.35%  "

# Row 43
$ws.Range("E43").Value = "  -0.68%  "

# Row 44
$ws.Range("E44").Value = "  -8.97%  "

# Row 45
$ws.Range("E45").Value = "  -2.90%  "

# Row 46
$ws.Range("D46").Value = "'3.40"
$ws.Range("E46").Value = "  +4.87%  "

# Row 47
$ws.Range("D47").Value = "'0.339"
$ws.Range("E47").Value = "  -0.95%  "

# Row 48
$ws.Range("D48").Value = "'0.0416"
$ws.Range("E48").Value = "  +0.20%  "

# Row 49
$ws.Range("D49").Value = "'0.129"
$ws.Range("E49").Value = "  -0.41%  "

# Row 50
$ws.Range("D50").Value = "'2.55"
$ws.Range("E50").Value = "  -2.28%  "

# Row 51
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.16%  "
